$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (Reinvestment *), shifting existing
# D:G columns to E:H.
$ws.Range("D1").EntireColumn.Insert()

# New column header
$ws.Range("D1").Value = "Cost Of Investment *"

# New column values ("Cost Of Investment *") for each distribution row
$ws.Range("D2").Value = 800000
$ws.Range("D3").Value = 1500000
$ws.Range("D4").Value = 2500000

# Match number formatting used by the other numeric columns (Gross / Reinvestment)
$ws.Range("D2:D4").NumberFormat = $ws.Range("E2:E4").NumberFormat

# Give the new column an explicit (non bestFit) width, matching column C
$ws.Range("D1").ColumnWidth = $ws.Range("C1").ColumnWidth

# Move the selection, matching the author's final cursor position
$ws.Range("D5").Select()
